# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell whose status
# was "Ready for handoff" is now "In Translation" (the shared string used
# by the Overview sheet's zh-cn/de-de columns and by the Status column on
# the zh-cn / de-de detail sheets). The source tool that builds this report
# also re-autosizes the Status-ish columns to fit the new (shorter) text,
# so their stored column widths shrink accordingly.

$wb = $excel.ActiveWorkbook

function Replace-CellText {
    param($ws, [string]$oldText, [string]$newText)

    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # Cast explicitly to string before comparing - PowerShell's -eq
            # coerces the right-hand operand to the left operand's type, so
            # comparing a Boolean cell value directly against a non-empty
            # string would otherwise evaluate to $true for every Boolean
            # cell (e.g. the "True" / "False" cells elsewhere in the sheet).
            $val = [string]$cell.Value2
            if ($val -ceq $oldText) {
                $cell.Value = $newText
            }
        }
    }
}

foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Replace-CellText $ws "Ready for handoff" "In Translation"
}

# Re-fit the columns that displayed the status text so they shrink to match
# the now-shorter "In Translation" label, same as the report generator does.
# (Excel's ColumnWidth setter quantizes to whole pixels, so 12.5 characters
# is the closest attainable width to the generator's fractional result.)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
